$p = $ppt.ActivePresentation
$s = $p.Slides.Item(29)
$shp = $s.Shapes.Item("Picture 5")
$shp.Top = 2645484 / 12700
